# arreglo error con funcion bCitas
# Adds a new "medico" user row to the Usuarios sheet (first sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New user record: Id=2, Usuario="medico", Password=123, Rol="medico"
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "medico"
$ws.Range("C3").Value = 123
$ws.Range("D3").Value = "medico"

# Match formatting of the row above (centered text) but without the border,
# by copying the existing row's format then stripping the border.
$ws.Range("A2:D2").Copy()
$ws.Range("A3:D3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A3:D3").Borders.LineStyle = -4142  # xlLineStyleNone

$excel.CutCopyMode = 0

$ws.Range("D3").Select()
